# Auto-generated edit script: applies value updates to the Kujata_Profits
# market-data workbook (columns H-N across sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
# as produced by the scheduled data-refresh runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$updates = @(
    @("H33", 439.2857),
    @("I33", 417.5),
    @("J33", 519.1667),
    @("K33", 417.5),
    @("L33", 519.1667),
    @("M33", -188.5),
    @("N33", -977.1667),
    @("H88", 1374798.4),
    @("I88", 753),
    @("K88", 753),
    @("M88", -347),
    @("H91", 1374798.4),
    @("I91", 753),
    @("K91", 753),
    @("M91", 651),
    @("H113", 3099.6667),
    @("I113", 3090),
    @("J113", 3104.5),
    @("K113", 3090),
    @("L113", 3104.5),
    @("M113", 164),
    @("N113", -9612.5)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$updates = @(
    @("H32", 4103.9683),
    @("I32", 3942.2712),
    @("J32", 6489),
    @("K32", 3942.2712),
    @("L32", 6489),
    @("M32", -3655.2712),
    @("N32", -7063),
    @("H45", 1121.6666),
    @("I45", 1069),
    @("K45", 1069),
    @("M45", -692),
    @("H74", 3148.75),
    @("I74", 2198.3333),
    @("J74", 6000),
    @("K74", 2198.3333),
    @("L74", 6000),
    @("M74", -1324.3333),
    @("N74", -7748),
    @("H77", 3148.75),
    @("I77", 2198.3333),
    @("J77", 6000),
    @("K77", 10991.6665),
    @("L77", 30000),
    @("M77", -6623.666499999999),
    @("N77", -38736),
    @("H122", 1548.0454),
    @("I122", 1207.8462),
    @("J122", 2039.4445),
    @("K122", 3623.5386),
    @("L122", 6118.333500000001),
    @("M122", -1173.5386),
    @("N122", -11018.3335)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$updates = @(
    @("H105", 47619960),
    @("I105", 47619960),
    @("K105", 47619960),
    @("M105", -47618213),
    @("H107", 1234.2667),
    @("I107", 800.1),
    @("J107", 2102.6),
    @("K107", 800.1),
    @("L107", 2102.6),
    @("M107", 1119.9),
    @("N107", -5942.6)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$updates = @(
    @("H7", 312.72726),
    @("I7", 129.6),
    @("K7", 129.6),
    @("M7", -16.59999999999999),
    @("H22", 140238.4),
    @("J22", 233600.33),
    @("L22", 233600.33),
    @("N22", -234300.33),
    @("H95", 10896.333),
    @("J95", 10896.333),
    @("L95", 10896.333),
    @("N95", -16388.333),
    @("H114", 23999.75),
    @("J114", 23999.75),
    @("L114", 23999.75),
    @("N114", -32677.75),
    @("H122", 717.64703),
    @("I122", 738.125),
    @("J122", 390),
    @("K122", 2214.375),
    @("L122", 1170),
    @("M122", 235.625),
    @("N122", -6070),
    @("H134", 18520440),
    @("I134", 2019.6086),
    @("J134", 125001350),
    @("K134", 6058.825800000001),
    @("L134", 375004050),
    @("M134", -3523.825800000001),
    @("N134", -375009120)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$updates = @(
    @("H74", 5453.6),
    @("J74", 5453.6),
    @("L74", 16360.8),
    @("N74", -18482.8),
    @("H77", 5453.6),
    @("J77", 5453.6),
    @("L77", 49082.4),
    @("N77", -59690.4),
    @("H81", 3700),
    @("I81", 0),
    @("K81", 0),
    @("H84", 3700),
    @("I84", 0),
    @("K84", 0),
    @("H122", 1048.3158),
    @("J122", 1059.3334),
    @("L122", 9534.000599999999),
    @("N122", -14434.0006)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$ws.Range("M81").ClearContents()
$ws.Range("M84").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$updates = @(
    @("H70", 40912220),
    @("I70", 41669420),
    @("J70", 40003576),
    @("K70", 41669420),
    @("L70", 40003576),
    @("M70", -41669150),
    @("N70", -40004116),
    @("H73", 40912220),
    @("I73", 41669420),
    @("J73", 40003576),
    @("K73", 41669420),
    @("L73", 40003576),
    @("M73", -41668484),
    @("N73", -40005448),
    @("H102", 1406.8462),
    @("I102", 1320.25),
    @("J102", 1545.4),
    @("K102", 1320.25),
    @("L102", 1545.4),
    @("M102", 301.75),
    @("N102", -4789.4),
    @("H122", 2000.125),
    @("I122", 2044.6875),
    @("J122", 1911),
    @("K122", 6134.0625),
    @("L122", 5733),
    @("M122", -3684.0625),
    @("N122", -10633),
    @("H132", 3317.1904),
    @("I132", 3369.5833),
    @("J132", 3247.3333),
    @("K132", 10108.7499),
    @("L132", 9741.999899999999),
    @("M132", -7578.749899999999),
    @("N132", -14801.9999),
    @("H134", 27536),
    @("J134", 27536),
    @("L134", 82608),
    @("N134", -87678)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$updates = @(
    @("H22", 1075),
    @("I22", 750),
    @("J22", 1400),
    @("K22", 750),
    @("L22", 1400),
    @("M22", -455),
    @("N22", -1990),
    @("H27", 1075),
    @("I27", 750),
    @("J27", 1400),
    @("K27", 750),
    @("L27", 1400),
    @("M27", -643),
    @("N27", -1614),
    @("H40", 3194.889),
    @("J40", 6481),
    @("L40", 6481),
    @("N40", -6753),
    @("H44", 12000),
    @("J44", 12000),
    @("L44", 12000),
    @("N44", -12912),
    @("H122", 50001840),
    @("I122", 62501700),
    @("J122", 2400),
    @("K122", 187505100),
    @("L122", 7200),
    @("M122", -187502650),
    @("N122", -12100)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$updates = @(
    @("H98", 19333.334),
    @("J98", 19333.334),
    @("L98", 19333.334),
    @("N98", -25323.334),
    @("H113", 312.5),
    @("I113", 220.66667),
    @("J113", 509.2857),
    @("K113", 662.00001),
    @("L113", 1527.8571),
    @("M113", 1507.99999),
    @("N113", -5867.8571),
    @("H132", 2389.4167),
    @("I132", 2087.65),
    @("J132", 3898.25),
    @("K132", 6262.950000000001),
    @("L132", 11694.75),
    @("M132", -3732.950000000001),
    @("N132", -16754.75)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

